$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B8: change from text "3" to a real number 3
$ws.Range("B8").Value = 3

# Add new row 9 with annotation data
$ws.Range("A9").Value = "Sunsi Wu"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "4"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "thank "
$ws.Range("D9").Value = "APC"
$ws.Range("E9").Value = "OTH"
$ws.Range("F9").Value = "41c93df3-3a59-4ce4-b94b-f420b7540586"
$ws.Range("G9").Value = "SJ19eUg0-_annotated.xlsx"
$ws.Range("H9").Value = "Thank the reviewer for the thoughtful feedback."
